$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at the top; existing rows 1-18 shift down to 2-19
$ws.Rows.Item(1).Insert()

# New row 1 only has a single cell A1 with the distance note
$ws.Range("A1").Value = "Distance from mirror to screen = 113 cm"

# Update selection to match the target state (A19 selected)
$ws.Range("A19").Select()
